$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4,  "paid",   148.95000000000002, 1),
    @(5,  "paid",   60.95,              9),
    @(6,  "booked", 101.95,             7),
    @(7,  "paid",   101.95,             7),
    @(8,  "booked", 37.97,              1),
    @(9,  "paid",   37.97,              1),
    @(10, "booked", 37.97,              1),
    @(11, "booked", 37.97,              1),
    @(12, "booked", 37.97,              1),
    @(13, "booked", 37.97,              1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $r
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
